$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.195.98"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.84%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.169.37"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.00%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.28"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.29%  "
$ws.Range("E6").Value = "  -1.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "70.07"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -4.54%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.579"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -5.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.18"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -7.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0925"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.98"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -4.46%  "
$ws.Range("E13").Value = "  -4.16%  "
$ws.Range("E14").Value = "  -2.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.492.98"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.88"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.804"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -4.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.166.65"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.011.00"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.89%  "
$ws.Range("E20").Value = "  -6.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.43"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.93"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.76"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -6.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "226.47"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.16%  "
$ws.Range("E25").Value = "  -6.54%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.89"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -5.17%  "
$ws.Range("E28").Value = "  -1.31%  "
$ws.Range("E29").Value = "  -2.32%  "
$ws.Range("E30").Value = "  +1.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "168.11"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.96"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.71"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +6.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0769"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.15"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -7.45%  "
$ws.Range("E36").Value = "  -3.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.102"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -8.29%  "
$ws.Range("E38").Value = "  -3.31%  "
$ws.Range("E39").Value = "  -5.30%  "
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.99"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -6.41%  "
$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.08"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.44"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "59.85"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -9.22%  "
$ws.Range("E44").Value = "  -4.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.30"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0976"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "97.77"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -5.77%  "
$ws.Range("E48").Value = "  -2.07%  "
$ws.Range("E49").Value = "  -2.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.21"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -8.37%  "
$ws.Range("E51").Value = "  -2.95%  "
